# Update "想去人数" (want-to-go count) values in the F column across sheets,
# reflecting refreshed counts pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 106
$ws1.Range("F5").Value  = 286
$ws1.Range("F7").Value  = 1148
$ws1.Range("F8").Value  = 400
$ws1.Range("F10").Value = 100
$ws1.Range("F12").Value = 24
$ws1.Range("F13").Value = 249
$ws1.Range("F15").Value = 138
$ws1.Range("F16").Value = 1311
$ws1.Range("F18").Value = 186
$ws1.Range("F19").Value = 297
$ws1.Range("F21").Value = 667
$ws1.Range("F22").Value = 1073
$ws1.Range("F23").Value = 55
$ws1.Range("F24").Value = 1936
$ws1.Range("F25").Value = 2480
$ws1.Range("F26").Value = 1253
$ws1.Range("F28").Value = 183
$ws1.Range("F29").Value = 354
$ws1.Range("F30").Value = 771
$ws1.Range("F31").Value = 754
$ws1.Range("F32").Value = 943
$ws1.Range("F36").Value = 347
$ws1.Range("F37").Value = 572
$ws1.Range("F38").Value = 706
$ws1.Range("F39").Value = 311
$ws1.Range("F40").Value = 204

# 演出 (Performances) sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 325

# 全部类型 (All types) sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 106
$ws4.Range("F8").Value  = 286
$ws4.Range("F12").Value = 1148
$ws4.Range("F13").Value = 400
$ws4.Range("F15").Value = 100
$ws4.Range("F17").Value = 249
$ws4.Range("F20").Value = 138
$ws4.Range("F21").Value = 1311
$ws4.Range("F23").Value = 186
$ws4.Range("F24").Value = 297
$ws4.Range("F26").Value = 1073
$ws4.Range("F27").Value = 2480
$ws4.Range("F29").Value = 1253
$ws4.Range("F34").Value = 183
$ws4.Range("F35").Value = 354
$ws4.Range("F36").Value = 771
$ws4.Range("F39").Value = 754
$ws4.Range("F40").Value = 943
$ws4.Range("F42").Value = 347
$ws4.Range("F43").Value = 572
$ws4.Range("F44").Value = 706
$ws4.Range("F45").Value = 311
$ws4.Range("F48").Value = 204
